$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.457.65"
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").Value = "3.025.56"
$ws.Range("E3").Value = "  +3.66%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "200.25"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").Value = "629.01"
$ws.Range("E6").Value = "  +4.88%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +0.29%  "

$ws.Range("D9").Value = "0.209"
$ws.Range("E9").Value = "  +4.48%  "

$ws.Range("D10").Value = "3.025.45"

$ws.Range("D11").Value = "0.436"

$ws.Range("E12").Value = "  -0.36%  "

$ws.Range("D13").Value = "5.13"
$ws.Range("E13").Value = "  +5.10%  "

$ws.Range("D14").Value = "3.580.29"
$ws.Range("E14").Value = "  +3.59%  "

$ws.Range("D15").Value = "29.21"
$ws.Range("E15").Value = "  +5.75%  "

$ws.Range("D16").Value = "76.325.12"
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").Value = "0.0000192"
$ws.Range("E17").Value = "  -0.36%  "

$ws.Range("D18").Value = "3.041.34"
$ws.Range("E18").Value = "  +4.29%  "

$ws.Range("D19").Value = "13.46"
$ws.Range("E19").Value = "  +3.26%  "

$ws.Range("E20").Value = "  +3.31%  "

$ws.Range("D21").Value = "373.60"
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("D22").Value = "4.36"
$ws.Range("E22").Value = "  +0.27%  "

$ws.Range("E23").Value = "  -1.44%  "

$ws.Range("D24").Value = "73.03"
$ws.Range("E24").Value = "  +2.43%  "

$ws.Range("D25").Value = "3.186.10"
$ws.Range("E25").Value = "  +3.84%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").Value = "4.38"
$ws.Range("E27").Value = "  +4.03%  "

$ws.Range("D28").Value = "9.89"
$ws.Range("E28").Value = "  +1.77%  "

$ws.Range("E29").Value = "  -0.63%  "

$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.18%  "

$ws.Range("D31").Value = "8.30"
$ws.Range("E31").Value = "  +7.29%  "

$ws.Range("D32").Value = "1.41"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").Value = "506.98"
$ws.Range("E33").Value = "  +0.59%  "

$ws.Range("E34").Value = "  +6.43%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").Value = "20.72"
$ws.Range("E36").Value = "  +2.38%  "

$ws.Range("D37").Value = "164.06"
$ws.Range("E37").Value = "  -0.90%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "0.385"
$ws.Range("E38").Value = "  +10.51%  "

$ws.Range("D39").Value = "20.00"
$ws.Range("E39").Value = "  +1.79%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "192.68"
$ws.Range("E40").Value = "  +6.27%  "

$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("E42").Value = "  -1.44%  "

$ws.Range("E43").Value = "  +0.38%  "

$ws.Range("D44").Value = "5.08"
$ws.Range("E44").Value = "  +1.72%  "

$ws.Range("D45").Value = "42.46"
$ws.Range("E45").Value = "  +5.66%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "1.67"
$ws.Range("E46").Value = "  +0.53%  "

$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").Value = "1.25"
$ws.Range("E47").Value = "  +5.39%  "

$ws.Range("D48").Value = "0.716"
$ws.Range("E48").Value = "  +8.49%  "

$ws.Range("D49").Value = "0.601"
$ws.Range("E49").Value = "  +4.72%  "

$ws.Range("D50").Value = "2.37"
$ws.Range("E50").Value = "  +1.36%  "

$ws.Range("D51").Value = "3.88"
$ws.Range("E51").Value = "  +4.26%  "
